$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "birds" -> "bird" label for pinguin's class
$ws.Range("B6").Value = "bird"

# New header for the "height" column (must be added to shared strings
# before "empty" so the shared-strings table order matches)
$ws.Range("E1").Value = "height"

# New "height" column values (column E, rows 2-8)
$ws.Range("E2").Value = 0.2
$ws.Range("E3").Value = 0.2
$ws.Range("E4").Value = 0.7
$ws.Range("E5").Value = 0.6
$ws.Range("E6").Value = 0.4
$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 0.1

# New header for the (empty) "empty" column
$ws.Range("D1").Value = "empty"

# Update the active selection to match the edited workbook
$ws.Range("D1").Select()
